# Rename sheets to reflect the switch from board-feet/cubic-feet units
# to thousand-board-feet / hundred-cubic-feet units.
$wb = $excel.ActiveWorkbook

$wsHarvest = $wb.Worksheets.Item("Harvest_BF")
$wsHarvest.Name = "Harvest_MBF"

$wsBfcf = $wb.Worksheets.Item("BFCF")
$wsBfcf.Name = "MBFCCF"

# Move the active/selected tab from HWP_MODEL_OPTIONS to MBFCCF, and update
# its selected cell to H30.
$wsBfcf.Select()
$wsBfcf.Range("H30").Select()
